$wb = $excel.ActiveWorkbook

# Sheet "NBR" - column C (Reaction_number) updated values for rows 2-20
$wsNBR = $wb.Worksheets.Item("NBR")
$nbrValues = @(836, 819, 745, 742, 739, 745, 735, 731, 735, 726, 714, 710, 716, 0, 704, 0, 0, 663, 655)
for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $row = $i + 2
    $wsNBR.Cells.Item($row, 3).Value = $nbrValues[$i]
}

# Sheet "BAR" - column C (Reaction_number) updated values for rows 2-20
$wsBAR = $wb.Worksheets.Item("BAR")
$barValues = @(691, 715, 751, 764, 760, 717, 712, 711, 705, 704, 706, 704, 699, 0, 697, 0, 0, 697, 696)
for ($i = 0; $i -lt $barValues.Length; $i++) {
    $row = $i + 2
    $wsBAR.Cells.Item($row, 3).Value = $barValues[$i]
}
